$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Rows 44-46 currently hold multi-run, tab-separated per-GC-cycle stats.
# They get collapsed down to the single summary numbers that used to live
# in rows 1-3 (93.83 / 47.71 / 773).
$t.Cell(44, 1).Range.Text = "93.83"
$t.Cell(45, 1).Range.Text = "47.71"
$t.Cell(46, 1).Range.Text = "773"

# Rows 1-3 (the old summary values) become "0M".
$t.Cell(1, 1).Range.Text = "0M"
$t.Cell(2, 1).Range.Text = "0M"
$t.Cell(3, 1).Range.Text = "0M"

# A handful of other standalone numeric updates.
$t.Cell(4, 1).Range.Text = "8420"
$t.Cell(7, 1).Range.Text = "0.04667"
$t.Cell(8, 1).Range.Text = "0.00277"
$t.Cell(12, 1).Range.Text = "47.71366"
